# 2021.04.17 17:48 hoon cli
# db스크립트에 단행본 insert 추가
#
# Sheet "단행본 후기" (Book review) gets a new second row holding the DB
# column names (matching the NOV_REVIEW/NOV_CONTENT style "insert script"
# header rows already used on the other table sheets of this workbook),
# pushing the existing sample data down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("단행본 후기")

# Insert a fresh row right under the Korean header row and push the
# existing 11 sample rows down.
$ws.Rows(2).Insert()

$ws.Range("A2").Value = "NOV_REVIEW_NUM"
$ws.Range("B2").Value = "NOV_NUM"
$ws.Range("C2").Value = "MEM_ID"
$ws.Range("D2").Value = "REVIEW_DATE"
$ws.Range("E2").Value = "NOV_CONTENT"

# Widen column A so the longer "NOV_REVIEW_NUM" header is fully visible.
$ws.Columns(1).ColumnWidth = 18.4

# Leave the cursor where the author left it after typing the new row.
$ws.Activate() | Out-Null
$ws.Range("C6").Select() | Out-Null

# Touch the page setup (paper size / orientation) like the source sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
